$d = $word.ActiveDocument

$replacements = @(
    @{old="49×14="; new="95×84="},
    @{old="49×36="; new="97×92="},
    @{old="29×12="; new="98×58="},
    @{old="89×77="; new="61×79="},
    @{old="75×42="; new="30×85="},
    @{old="59×53="; new="25×29="},
    @{old="54×92="; new="20×99="},
    @{old="16×38="; new="67×45="},
    @{old="33×94="; new="32×28="},
    @{old="29×46="; new="24×84="},
    @{old="15×27="; new="39×90="},
    @{old="12×13="; new="41×36="},
    @{old="36×39="; new="86×86="},
    @{old="25×68="; new="74×85="},
    @{old="27×73="; new="86×68="},
    @{old="45×91="; new="63×47="},
    @{old="36×85="; new="13×68="},
    @{old="67×25="; new="91×46="},
    @{old="79×71="; new="89×50="},
    @{old="32×45="; new="61×76="},
    @{old="52×73="; new="74×73="},
    @{old="96×78="; new="57×48="},
    @{old="15×53="; new="70×22="},
    @{old="64×30="; new="73×34="},
    @{old="52×61="; new="27×29="}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
